# Added filtering options for the Component Analysis
# Clear out cells that fall outside the filtered range for the
# Component Analysis on rows 2, 3, 5, 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5:K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
